$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 3-8 (the data rows are reordered chronologically by
# the "Fecha" column D, carrying their related M, N, O, P, Q, S values along).
$data = @(
    @{ Row = 3; D = 44162; M = 120; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 },
    @{ Row = 4; D = 44176; M = 250; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 },
    @{ Row = 5; D = 44208; M = 210; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 },
    @{ Row = 6; D = 44309; M = 300; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 },
    @{ Row = 7; D = 44351; M = 300; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 },
    @{ Row = 8; D = 44397; M = 60;  N = 11000; O = 11000; P = 11000; Q = "`$/caja 14 kilos";           S = 786 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value  = $item.D   # D - Fecha
    $ws.Cells.Item($r, 13).Value = $item.M   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $item.N   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $item.O   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $item.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $item.Q   # Q - Unidad de comercializacion
    $ws.Cells.Item($r, 19).Value = $item.S   # S - Precio $/Kg
}
